$d = $word.ActiveDocument

# 1. Append a new highlighted run "(in progress)" to the paragraph
#    "Implement initial ideas for the basis of each level design individually."
$target = $d.Paragraphs.Item(6)
$r = $target.Range
$insertPos = $r.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("(in progress)")

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Highlight = 1
[void]$find.Execute("(in progress)", $true, $false, $false, $false, $false,
                     $true, 1, $false, "(in progress)", 2)

# 2. Apply strikethrough formatting (paragraph mark + run) to the three
#    completed milestone paragraphs.
$strikeParagraphs = @(
    "Program ability ball mechanics to behave as intended.",
    "Program the success state to occur when the player reaches the goal region with the objective ball.",
    "Program a fail state to occur when the player does not complete a stage within the given limitations."
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($needle in $strikeParagraphs) {
        if ($text -eq $needle) {
            $p.Range.Font.StrikeThrough = 1
        }
    }
}
